$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a text value on a cell that looks like a percentage,
# preventing Excel's auto-conversion to a numeric percent value while
# preserving the cell's original style (border formatting).
function Set-PercentText($cell, $value, $donor) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
    $cell.Borders.LineStyle = $donor.Borders.LineStyle
}

# --- Plain text/date value updates ---
$ws.Range("E2").Value = "2026-02-26 20:18:20"
$ws.Range("E3").Value = "2026-02-26 20:18:23"
$ws.Range("L3").Value = "16.9 km/h - 114º 19:47 TU"
$ws.Range("E4").Value = "2026-02-26 20:18:25"
$ws.Range("O4").Value = "10.9 °C"
$ws.Range("E5").Value = "2026-02-26 20:18:28"
$ws.Range("E6").Value = "2026-02-26 20:18:30"
$ws.Range("E7").Value = "2026-02-26 20:18:32"
$ws.Range("E8").Value = "2026-02-26 20:18:35"
$ws.Range("E9").Value = "2026-02-26 20:18:37"
$ws.Range("O9").Value = "12.3 °C"
$ws.Range("E10").Value = "2026-02-26 20:18:38"
$ws.Range("O10").Value = "9.7 °C"
$ws.Range("E11").Value = "2026-02-26 20:18:39"
$ws.Range("E12").Value = "2026-02-26 20:18:40"
$ws.Range("E13").Value = "2026-02-26 20:18:41"
$ws.Range("E14").Value = "2026-02-26 20:18:42"
$ws.Range("N14").Value = "7.0 °C 19:35 TU"
$ws.Range("O14").Value = "11.6 °C"
$ws.Range("E15").Value = "2026-02-26 20:18:44"
$ws.Range("N15").Value = "6.6 °C 19:46 TU"
$ws.Range("O15").Value = "11.7 °C"
$ws.Range("E16").Value = "2026-02-26 20:18:45"
$ws.Range("E17").Value = "2026-02-26 20:18:46"
$ws.Range("E18").Value = "2026-02-26 20:18:47"
$ws.Range("O18").Value = "12.0 °C"
$ws.Range("E19").Value = "2026-02-26 20:18:48"
$ws.Range("E20").Value = "2026-02-26 20:18:49"
$ws.Range("O20").Value = "2.7 °C"
$ws.Range("E21").Value = "2026-02-26 20:18:52"
$ws.Range("E22").Value = "2026-02-26 20:18:54"
$ws.Range("K22").Value = "17.5 MJ/m2"
$ws.Range("O22").Value = "2.3 °C"
$ws.Range("E23").Value = "2026-02-26 20:18:56"
$ws.Range("E24").Value = "2026-02-26 20:18:59"
$ws.Range("O24").Value = "10.6 °C"
$ws.Range("E25").Value = "2026-02-26 20:19:01"
$ws.Range("O25").Value = "5.3 °C"
$ws.Range("E26").Value = "2026-02-26 20:19:04"
$ws.Range("E27").Value = "2026-02-26 20:19:06"
$ws.Range("O27").Value = "5.2 °C"
$ws.Range("E28").Value = "2026-02-26 20:19:09"
$ws.Range("O28").Value = "11.0 °C"
$ws.Range("E29").Value = "2026-02-26 20:19:11"
$ws.Range("N29").Value = "8.3 °C 19:56 TU"
$ws.Range("E30").Value = "2026-02-26 20:19:13"
$ws.Range("O30").Value = "12.2 °C"
$ws.Range("E31").Value = "2026-02-26 20:19:16"
$ws.Range("O31").Value = "11.9 °C"
$ws.Range("E32").Value = "2026-02-26 20:19:19"
$ws.Range("O32").Value = "8.2 °C"
$ws.Range("E33").Value = "2026-02-26 20:19:21"
$ws.Range("E34").Value = "2026-02-26 20:19:24"
$ws.Range("O34").Value = "4.9 °C"
$ws.Range("E35").Value = "2026-02-26 20:19:26"
$ws.Range("E36").Value = "2026-02-26 20:19:28"
$ws.Range("N36").Value = "9.2 °C 19:38 TU"
$ws.Range("E37").Value = "2026-02-26 20:19:31"
$ws.Range("E38").Value = "2026-02-26 20:19:33"
$ws.Range("O38").Value = "11.3 °C"
$ws.Range("E39").Value = "2026-02-26 20:19:35"
$ws.Range("O39").Value = "2.9 °C"
$ws.Range("E40").Value = "2026-02-26 20:19:38"
$ws.Range("E41").Value = "2026-02-26 20:19:40"
$ws.Range("O41").Value = "11.2 °C"
$ws.Range("E42").Value = "2026-02-26 20:19:42"
$ws.Range("E43").Value = "2026-02-26 20:19:45"
$ws.Range("E44").Value = "2026-02-26 20:19:47"
$ws.Range("E45").Value = "2026-02-26 20:19:50"
$ws.Range("E46").Value = "2026-02-26 20:19:52"

# --- Percentage text value updates (special handling) ---
Set-PercentText $ws.Range("H6") "83%" $ws.Range("I6")
Set-PercentText $ws.Range("H11") "68%" $ws.Range("I11")
Set-PercentText $ws.Range("H14") "88%" $ws.Range("I14")
Set-PercentText $ws.Range("H19") "46%" $ws.Range("I19")
Set-PercentText $ws.Range("H26") "41%" $ws.Range("I26")
Set-PercentText $ws.Range("H27") "41%" $ws.Range("I27")
Set-PercentText $ws.Range("H28") "78%" $ws.Range("I28")
Set-PercentText $ws.Range("H30") "86%" $ws.Range("I30")
Set-PercentText $ws.Range("H32") "65%" $ws.Range("I32")
Set-PercentText $ws.Range("H37") "73%" $ws.Range("I37")
Set-PercentText $ws.Range("H44") "53%" $ws.Range("I44")
